$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (REG VAL)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.050949134"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.14144452"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 76
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 2

# Row 5 (REG TEST)
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0.08170867"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "0.22512843"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").Value = 21
$ws.Range("E5").Value = 9
